$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("working")

# Update D6 value from 100 to 1000
$ws.Range("D6").Value = 1000

# Update the active selection to F4 (matches the post-edit cursor position)
$ws.Range("F4").Select()
